$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.818.92'
$ws.Range("E2").Value = '  +0.05%  '

$ws.Range("D3").Value = '2.313.46'
$ws.Range("E3").Value = '  +0.67%  '

$ws.Range("E4").Value = '  -0.22%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.37'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.43%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '107.50'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.62%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.629'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.17%  '

$ws.Range("E8").Value = '  -0.17%  '

$ws.Range("E9").Value = '  +0.68%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.17'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.42%  '

$ws.Range("E11").Value = '  +0.52%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.41'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.45%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.108'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.69%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.992'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.04%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.27'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.82%  '

$ws.Range("D16").Value = '2.664.57'
$ws.Range("E16").Value = '  +0.65%  '

$ws.Range("D17").Value = '2.314.37'
$ws.Range("E17").Value = '  -0.31%  '

$ws.Range("D18").Value = '42.794.04'
$ws.Range("E18").Value = '  +0.13%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.49'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.77%  '

$ws.Range("E20").Value = '  -0.22%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.11'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -13.61%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.68'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.68%  '

$ws.Range("E23").Value = '  -1.62%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '266.65'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.04%  '

$ws.Range("E25").Value = '  +0.65%  '

$ws.Range("E26").Value = '  -0.04%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.73'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +11.63%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.07'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.21%  '

$ws.Range("E29").Value = '  -2.28%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.56'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.95%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.49'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.55%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '166.41'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.03%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0878'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.12%  '

$ws.Range("E34").Value = '  +4.71%  '

$ws.Range("E35").Value = '  -0.94%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.73'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.56%  '

$ws.Range("E37").Value = '  -1.36%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.82'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.84%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.67'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.87%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.62'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.46%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '104.08'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +8.88%  '

$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.234'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.79%  '

$ws.Range("B44").Value = 'MultiversX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '70.91'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.52%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.07'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.66%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.40%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '113.03'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.77%  '

$ws.Range("D48").Value = '1.662.52'
$ws.Range("E48").Value = '  -2.85%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '77.28'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.82%  '

$ws.Range("E50").Value = '  +3.01%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.81'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.23%  '
